$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '49.824.33'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.61%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.630.59'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.10%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.18%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.47%  '

$ws.Range('E7').Value = '  +1.38%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.557'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.67%  '

$ws.Range('E10').Value = '  +2.70%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.57'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0818'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.77%  '

$ws.Range('E13').Value = '  +0.74%  '

$ws.Range('E14').Value = '  +2.51%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.046.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.24%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.626.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.02%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.879'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '49.770.52'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.82%  '

$ws.Range('E19').Value = '  +11.04%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.95%  '

$ws.Range('E21').Value = '  +1.53%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0959'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.22%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '281.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.42%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.08%  '

$ws.Range('E25').Value = '  +2.38%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.64'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.74%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.15%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.92%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.92'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.98%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.49'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.89%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.144'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.17%  '

$ws.Range('E32').Value = '  +0.38%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.86%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.33%  '

$ws.Range('E35').Value = '  +0.05%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0794'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.56%  '

$ws.Range('E37').Value = '  +5.90%  '

$ws.Range('E38').Value = '  +2.48%  '

$ws.Range('E39').Value = '  +7.67%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '123.34'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.12%  '

$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.22'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.52%  '

$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.01%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0314'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.063.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.78%  '

$ws.Range('E47').Value = '  +15.39%  '

$ws.Range('E48').Value = '  +7.53%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.70%  '

$ws.Range('E50').Value = '  +4.06%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.50'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.57%  '
